$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme
try {
    $tcs.Save("/tmp/work/saved_colors.xml")
    Write-Output "Save OK"
} catch {
    Write-Output "ERR: $_"
}
